$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(336).Insert()

$ws.Cells.Item(336, 1).Value = 9
$ws.Cells.Item(336, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(336, 3).Value = "Metropolitana"
$ws.Cells.Item(336, 4).Value = 44809
$ws.Cells.Item(336, 5).Value = 13
$ws.Cells.Item(336, 6).Value = 100112039
$ws.Cells.Item(336, 7).Value = "Ciboulette"
$ws.Cells.Item(336, 8).Value = "Sin especificar"
$ws.Cells.Item(336, 9).Value = "Primera"
$ws.Cells.Item(336, 10).Value = 250
$ws.Cells.Item(336, 11).Value = 1000
$ws.Cells.Item(336, 12).Value = 1200
$ws.Cells.Item(336, 13).Value = 1100
$ws.Cells.Item(336, 14).Value = "$/docena de atados"
$ws.Cells.Item(336, 15).Value = "Región Metropolitana"
$ws.Cells.Item(336, 16).Value = 367
$ws.Cells.Item(336, 17).Value = 3
$ws.Cells.Item(336, 18).Value = "Hortaliza"
